$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.199.82'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.561.00'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '210.79'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').Value = '0.489'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '22.05'
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').Value = '1.782.98'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '1.560.46'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').Value = '3.76'
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').Value = '0.516'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '27.170.26'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').Value = '61.94'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('D18').Value = '7.46'
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('D19').Value = '216.76'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '0.0₃0700'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').Value = '9.23'
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('D24').Value = '1.93'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').Value = '152.26'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').Value = '6.62'
$ws.Range('E26').Value = '  +0.47%  '
$ws.Range('D27').Value = '15.05'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +1.94%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D31').Value = '0.0468'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.450.93'
$ws.Range('E33').Value = '  +2.18%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '3.17'
$ws.Range('E34').Value = '  +1.74%  '
$ws.Range('E35').Value = '  +4.40%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Value = '2.33'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').Value = '0.544'
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '0.809'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.33'
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').Value = '0.992'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').Value = '64.10'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').Value = '1.73'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '1.696.66'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '85.65'
$ws.Range('E48').Value = '  -1.51%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0526'
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('D51').Value = '0.0951'
$ws.Range('E51').Value = '  -0.56%  '
